# "Is hashed" column (G) was converted from native boolean cells to
# text cells holding the literal words TRUE / FALSE (as shared strings),
# keeping the same logical meaning (previously boolean TRUE -> text "TRUE",
# boolean FALSE -> text "FALSE").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$trueRows  = @(2,3,4,5,6,7,8,9,10,11,14,15,16)
$falseRows = @(12,13,17,18,19,20)

foreach ($r in $trueRows) {
    $cell = $ws.Range("G$r")
    $cell.Formula = "=""TRUE"""
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

foreach ($r in $falseRows) {
    $cell = $ws.Range("G$r")
    $cell.Formula = "=""FALSE"""
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}
